$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$apos = "'"               # leading apostrophe forces Excel to keep text as-is
$sub3 = [string][char]0x2083   # subscript three, U+2083 (cast to string to avoid numeric coercion)

$ws.Range("D2").Value = $apos + '67.279.00'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.19%  '

$ws.Range("D3").Value = $apos + '3.108.96'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.40%  '

$ws.Range("D4").Value = $apos + '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = $apos + '574.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.90%  '

$ws.Range("D6").Value = $apos + '178.73'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.04%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").Value = $apos + '3.107.74'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.47%  '

$ws.Range("E9").Value = '  -0.90%  '

$ws.Range("E10").Value = '  -0.92%  '

$ws.Range("E11").Value = '  +0.32%  '

$ws.Range("E12").Value = '  -1.51%  '

$ws.Range("E13").Value = '  -1.50%  '

$ws.Range("D14").Value = $apos + '36.25'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.74%  '

$ws.Range("E15").Value = '  +0.25%  '

$ws.Range("D16").Value = $apos + '3.626.18'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.45%  '

$ws.Range("D17").Value = $apos + '67.190.23'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.12%  '

$ws.Range("D18").Value = $apos + '7.05'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.31%  '

$ws.Range("D19").Value = $apos + '3.107.88'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.47%  '

$ws.Range("D20").Value = $apos + '16.63'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.39%  '

$ws.Range("D21").Value = $apos + '491.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.28%  '

$ws.Range("D22").Value = $apos + '7.75'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.16%  '

$ws.Range("D23").Value = $apos + '0.688'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.14%  '

$ws.Range("D24").Value = $apos + '83.61'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.07%  '

$ws.Range("B25").Value = 'Fetch.AI'
$ws.Range("C25").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D25").Value = $apos + '2.27'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.35%  '

$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").Value = $apos + '12.62'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.75%  '

$ws.Range("D27").Value = $apos + '10.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.38%  '

$ws.Range("E28").Value = '  +0.03%  '

$ws.Range("D29").Value = $apos + '7.98'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.66%  '

$ws.Range("E30").Value = '  -0.21%  '

$ws.Range("E31").Value = '  -2.15%  '

$ws.Range("D32").Value = $apos + '28.19'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.04%  '

$ws.Range("D33").Value = $apos + '0.112'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.34%  '

$ws.Range("D34").Value = $apos + '0.0' + $sub3 + '0945'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.01%  '

$ws.Range("D35").Value = $apos + '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.01%  '

$ws.Range("D36").Value = $apos + '47.53'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.26%  '

$ws.Range("E37").Value = '  -1.81%  '

$ws.Range("E38").Value = '  -2.91%  '

$ws.Range("E39").Value = '  +2.65%  '

$ws.Range("E40").Value = '  +0.35%  '

$ws.Range("D41").Value = $apos + '49.16'
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").Value = $apos + '0.124'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.32%  '

$ws.Range("B43").Value = 'Cosmos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D43").Value = $apos + '8.30'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.57%  '

$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").Value = $apos + '2.73'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.61%  '

$ws.Range("D45").Value = $apos + '2.802.90'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.45%  '

$ws.Range("D46").Value = $apos + '372.23'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.33%  '

$ws.Range("E47").Value = '  -0.82%  '

$ws.Range("D48").Value = $apos + '135.77'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.55%  '

$ws.Range("D50").Value = $apos + '25.62'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.34%  '

$ws.Range("D51").Value = $apos + '2.30'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.32%  '
